$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns D (Diferenca) and E (porcentagem), matching the
# existing bold/bordered header style used by A1:C1
$ws.Range("A1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Diferenca"
$ws.Range("E1").Value = "porcentagem"

# Refresh existing rows (Priori counts changed) and append new rows for
# the seats/"cadeiras" that were added, plus the new Diferenca/porcentagem columns
$ws.Range("A2").Value = "CAFI"
$ws.Range("B2").Value = 9884
$ws.Range("C2").Value = 8354
$ws.Range("D2").Value = 1530
$ws.Range("E2").Value = 84.52043707001215

$ws.Range("A3").Value = "CBAR"
$ws.Range("B3").Value = 20476
$ws.Range("C3").Value = 8234
$ws.Range("D3").Value = 12242
$ws.Range("E3").Value = 40.21293221332292

$ws.Range("A4").Value = "CBJA"
$ws.Range("B4").Value = 21167
$ws.Range("C4").Value = 7211
$ws.Range("D4").Value = 13956
$ws.Range("E4").Value = 34.06718004440875

$ws.Range("A5").Value = "CGAR"
$ws.Range("B5").Value = 11963
$ws.Range("C5").Value = 8866
$ws.Range("D5").Value = 3097
$ws.Range("E5").Value = 74.11184485496949

$ws.Range("A6").Value = "CPES"
$ws.Range("B6").Value = 13878
$ws.Range("C6").Value = 9519
$ws.Range("D6").Value = 4359
$ws.Range("E6").Value = 68.59057501080848

$ws.Range("A7").Value = "CPLT"
$ws.Range("B7").Value = 7212
$ws.Range("C7").Value = 4620
$ws.Range("D7").Value = 2592
$ws.Range("E7").Value = 64.05990016638935

$ws.Range("A8").Value = "CABL"
$ws.Range("B8").Value = 3765
$ws.Range("C8").Value = 2905
$ws.Range("D8").Value = 860
$ws.Range("E8").Value = 77.15803452855245

$ws.Range("A9").Value = "CCAR"
$ws.Range("B9").Value = 12203
$ws.Range("C9").Value = 4713
$ws.Range("D9").Value = 7490
$ws.Range("E9").Value = 38.62165041383266

$ws.Range("A10").Value = "CCSA"
$ws.Range("B10").Value = 9977
$ws.Range("C10").Value = 7472
$ws.Range("D10").Value = 2505
$ws.Range("E10").Value = 74.89225218001403

$ws.Range("A11").Value = "CIGR"
$ws.Range("B11").Value = 5045
$ws.Range("C11").Value = 3374
$ws.Range("D11").Value = 1671
$ws.Range("E11").Value = 66.8780971258672

$ws.Range("A12").Value = "CIPJ"
$ws.Range("B12").Value = 15033
$ws.Range("C12").Value = 8612
$ws.Range("D12").Value = 6421
$ws.Range("E12").Value = 57.28730127053815

$ws.Range("A13").Value = "CJBG"
$ws.Range("B13").Value = 3039
$ws.Range("C13").Value = 2641
$ws.Range("D13").Value = 398
$ws.Range("E13").Value = 86.90358670615333

$ws.Range("A14").Value = "COLI"
$ws.Range("B14").Value = 3832
$ws.Range("C14").Value = 2992
$ws.Range("D14").Value = 840
$ws.Range("E14").Value = 78.07933194154488

$ws.Range("A15").Value = "CPMR"
$ws.Range("B15").Value = 3029
$ws.Range("C15").Value = 1935
$ws.Range("D15").Value = 1094
$ws.Range("E15").Value = 63.88246946186861

$ws.Range("A16").Value = "CREC"
$ws.Range("B16").Value = 55364
$ws.Range("C16").Value = 34587
$ws.Range("D16").Value = 20777
$ws.Range("E16").Value = 62.47200346795751

$ws.Range("A17").Value = "CVSA"
$ws.Range("B17").Value = 22583
$ws.Range("C17").Value = 15069
$ws.Range("D17").Value = 7514
$ws.Range("E17").Value = 66.72718416507992

$ws.Range("A18").Value = "REIF"
$ws.Range("B18").Value = 19344
$ws.Range("C18").Value = 7482
$ws.Range("D18").Value = 11862
$ws.Range("E18").Value = 38.67866004962779

